# Mistake in manual gating (Freq of parent instead of count)
#
# Sheet1 mistakenly reported "Freq of parent" percentages in the Pop7
# column (column G) instead of the raw event Count. Fix: duplicate the
# sheet as "Sheet2" (placed after "Sheet1") with the corrected Count
# values in column G, and leave Sheet1's data untouched (only its
# selection/active-tab state changes because Sheet2 becomes the active
# sheet).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the corrected sheet right after Sheet1; Excel auto-names it "Sheet2".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Duplicate all of Sheet1's data (headers + values) onto the new sheet.
$ws1.UsedRange.Copy($ws2.Range("A1"))

# Correct column G (Pop7) with the true Count values instead of the
# mistaken Freq of parent percentages.
$correctedCounts = @(340, 229, 132, 505, 274, 746, 160, 172, 19, 103, 170, 142, 183, 441, 117)
for ($i = 0; $i -lt $correctedCounts.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 7).Value = $correctedCounts[$i]
}

# Restore each sheet's on-screen selection: Sheet1 is no longer the
# active tab, Sheet2 (the fixed data) becomes active.
$ws1.Range("C31").Select() | Out-Null
$ws2.Range("D28").Select() | Out-Null
